$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / first column relabeling (type, vi_mod, vi_sev) with bold styling ---
$ws.Range("A1").Value = "type"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Value = "vi_mod"
$ws.Range("A2").Font.Bold = $true

$ws.Range("A3").Value = "vi_sev"
$ws.Range("A3").Font.Bold = $true

# --- Row 2 data refresh (tornado diagram source values) ---
$ws.Range("B2").Value = 102.89816838
$ws.Range("C2").Value = 302.17330000000004
$ws.Range("D2").Value = 448.85769599999998
$ws.Range("E2").Value = 86.004435599999994
$ws.Range("F2").Value = 103.16822442029093
$ws.Range("G2").Value = 117.75637962589533
$ws.Range("H2").Value = 1673.2624000000001
$ws.Range("I2").Value = 1156.2883199999999
$ws.Range("J2").Value = 195.33871899624728
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 15.66
$ws.Range("M2").Value = 24.189983999999999
$ws.Range("N2").Value = 1.951962
$ws.Range("O2").Value = 587.42976319999991
$ws.Range("P2").Value = 76.716332800000004
$ws.Range("Q2").Value = 475.56818399999997
$ws.Range("R2").Value = 152.5275584
$ws.Range("S2").Value = 7630.3786967924143
$ws.Range("T2").Value = 0

# --- Column width for new tornado-diagram label column ---
$ws.Range("S1").ColumnWidth = 19.1667

# --- Selection / view state ---
$ws.Range("V1").Select()
